# Auto-generated script to update cryptos.xlsx price/volume columns
# per commit "Updated cryptos list on Wed Oct 23 14:28:41 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.681.07"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "2.575.20"
$ws.Range("E3").Value = "  -1.98%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E4").Value = "  -0.18%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "584.70"
$c.ClearFormats()
$ws.Range("E5").Value = "  -1.44%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "168.39"
$c.ClearFormats()
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D9").Value = "2.572.60"
$ws.Range("E9").Value = "  -2.00%  "
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("E11").Value = "  +0.30%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.356"
$c.ClearFormats()
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("E13").Value = "  -1.02%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "26.84"
$c.ClearFormats()
$ws.Range("E14").Value = "  -2.74%  "
$ws.Range("E16").Value = "  -1.44%  "
$ws.Range("D17").Value = "66.469.34"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "2.585.62"
$ws.Range("E18").Value = "  -1.85%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.42"
$c.ClearFormats()
$ws.Range("E19").Value = "  -6.00%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.73"
$c.ClearFormats()
$ws.Range("E20").Value = "  -3.57%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "350.81"
$c.ClearFormats()
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("E22").Value = "  -1.68%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "4.61"
$c.ClearFormats()
$ws.Range("E23").Value = "  -1.09%  "
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("E26").Value = "  -0.98%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.90"
$c.ClearFormats()
$ws.Range("E27").Value = "  -8.84%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("D30").Value = "0.0₃0991"
$ws.Range("E30").Value = "  -1.40%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "531.84"
$c.ClearFormats()
$ws.Range("E31").Value = "  -3.10%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "8.19"
$c.ClearFormats()
$ws.Range("E32").Value = "  +3.60%  "
$ws.Range("E33").Value = "  -1.96%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.84"
$c.ClearFormats()
$ws.Range("E34").Value = "  -3.12%  "
$ws.Range("E35").Value = "  -2.79%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("E37").Value = "  -2.43%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "156.74"
$c.ClearFormats()
$ws.Range("E38").Value = "  +0.06%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "18.80"
$c.ClearFormats()
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("E40").Value = "  -1.72%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "18.33"
$c.ClearFormats()
$ws.Range("E41").Value = "  +2.14%  "
$ws.Range("E42").Value = "  -0.23%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.78"
$c.ClearFormats()
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("E45").Value = "  +0.88%  "
$ws.Range("E46").Value = "  -3.32%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "149.26"
$c.ClearFormats()
$ws.Range("E47").Value = "  -1.39%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.567"
$c.ClearFormats()
$ws.Range("E48").Value = "  -2.12%  "
$ws.Range("E49").Value = "  -1.37%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.73"
$c.ClearFormats()
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("E51").Value = "  -0.93%  "
